# Katalog.xlsx edit: add a "Konsolen" (consoles) pricing line to the
# Draht_Matten sheet, update the Montageart options and the Gesamtpreis
# formula text accordingly, and leave the workbook with Draht_Matten as
# the active / selected sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Draht_Matten")

# The "Montageart" options (row 5) no longer offer "Dübelplatte" but a
# "Konsole" option instead.
$ws.Range("D5").Value = "Betonieren:0, Konsole:1"

# Insert a brand-new row 6 for the "Konsolen" variable, pushing the
# existing "Montage (€/m)" / "Gesamtpreis" rows down by one.
$ws.Rows.Item(6).Insert()
$ws.Range("A6").Value = "Auswahl"
$ws.Range("B6").Value = "Konsolen"
$ws.Range("C6").Value = "P_Konsolen"
$ws.Range("D6").Value = "Schwer:50, leicht:30"

# The total-price formula (documentation text, row 8 after the insert)
# now factors the console price in.
$ws.Range("E8").Value = "(L * P_Matte) + ((math.ceil(L/2.5)+1) * (P_Saeule + (P_Fund*P_Konsole))) + (L * P_Arbeit)"

# Widen the "Variable" / "Optionen" columns so the new, longer entries
# are readable.
$ws.Columns.Item(3).ColumnWidth = 14.3
$ws.Columns.Item(4).ColumnWidth = 13.14

# Leave Draht_Matten as the selected / active sheet with D12 selected.
$ws.Activate()
$ws.Range("D12").Select()
